# Applies the vocabulary.xlsx update described by the commit
# "new .ttl from Google sheet has been generated":
#   - bump dct:modified timestamp (B19)
#   - append vocab:1178 to the skos:broader lists of a few rows (D26, D27, D78)
#   - lower-case a handful of WHO GHO Observatory Framework labels (rows 176-186)
#     and remap their skos:broader from vocab:1152 -> vocab:1154
#   - rename/repurpose row 187 from "No research category assigned" to "disease"
#     and drop its skos:broader value
#   - append 18 new terms (rows 188-205) describing diseases/pathogens

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- dct:modified^^xsd:datetime ---------------------------------------
$ws.Range("B19").Value = "2023-08-22T13:09:45+00:00"

# --- append vocab:1178 to existing skos:broader lists -----------------
$ws.Range("D26").Value = "vocab:1000,vocab:1007,vocab:1020,vocab:1024,vocab:1037,vocab:1050,vocab:1058,vocab:1078,vocab:1147,vocab:1148,vocab:1154,vocab:1178"
$ws.Range("D27").Value = "vocab:1000,vocab:1007,vocab:1020,vocab:1024,vocab:1037,vocab:1050,vocab:1058,vocab:1078,vocab:1089,vocab:1099,vocab:1100,vocab:1154,vocab:1178"
$ws.Range("D78").Value = "vocab:1050,vocab:1058,vocab:1089,vocab:1166,vocab:1178"

# --- rows 176-186: lower-case labels + broader -> vocab:1154 ----------
$ws.Range("B176").Value = "cause"
$ws.Range("D176").Value = "vocab:1154"

$ws.Range("B177").Value = "solution"
$ws.Range("D177").Value = "vocab:1154"

$ws.Range("B178").Value = "drug development"
$ws.Range("D178").Value = "vocab:1154"

$ws.Range("B179").Value = "diagnostics/screening"
$ws.Range("D179").Value = "vocab:1154"

$ws.Range("B180").Value = "device"
$ws.Range("D180").Value = "vocab:1154"

$ws.Range("B181").Value = "vaccine"
$ws.Range("D181").Value = "vocab:1154"

$ws.Range("B182").Value = "vector control"
$ws.Range("D182").Value = "vocab:1154"

$ws.Range("B183").Value = "basic Science"
$ws.Range("D183").Value = "vocab:1154"

$ws.Range("B184").Value = "implementation"
$ws.Range("D184").Value = "vocab:1154"

$ws.Range("B185").Value = "evaluation"
$ws.Range("D185").Value = "vocab:1154"

$ws.Range("B186").Value = "no research category assigned"
$ws.Range("D186").Value = "vocab:1154"

# --- row 187: repurposed term, no skos:broader -------------------------
$ws.Range("B187").Value = "disease"
$ws.Range("D187").Value = ""

# --- new rows 188-205 ---------------------------------------------------
$newRows = @(
    @{ Row = 188; A = "vocab:1167"; B = "COVID-19"; D = "vocab:1166"; F = "http://purl.org/zonmw/covid19/10270" },
    @{ Row = 189; A = "vocab:1168"; B = "Crimean-Congo haemorrhagic fever"; D = "vocab:1166" },
    @{ Row = 190; A = "vocab:1169"; B = "Ebola virus disease"; D = "vocab:1166" },
    @{ Row = 191; A = "vocab:1170"; B = "Marburg virus disease"; D = "vocab:1166" },
    @{ Row = 192; A = "vocab:1171"; B = "Lassa fever"; D = "vocab:1166" },
    @{ Row = 193; A = "vocab:1172"; B = "Middle East Respiratory Syndrome Coronavirus (MERS-CoV)"; D = "vocab:1166" },
    @{ Row = 194; A = "vocab:1173"; B = "Severe Acute Respiratory Syndrome (SARS)"; D = "vocab:1166" },
    @{ Row = 195; A = "vocab:1174"; B = "Nipah and henipaviral disease"; D = "vocab:1166" },
    @{ Row = 196; A = "vocab:1175"; B = "Rift Valley Fever"; D = "vocab:1166" },
    @{ Row = 197; A = "vocab:1176"; B = "Zika"; D = "vocab:1166" },
    @{ Row = 198; A = "vocab:1177"; B = "Congenital Zika virus disease"; D = "vocab:1166" },
    @{ Row = 199; A = "vocab:1178"; B = "pathogens" },
    @{ Row = 200; A = "vocab:1179"; B = "Coronavirus"; D = "vocab:1178" },
    @{ Row = 201; A = "vocab:1180"; B = "SARS-CoV"; D = "vocab:1178" },
    @{ Row = 202; A = "vocab:1181"; B = "SARSr-CoV"; D = "vocab:1178" },
    @{ Row = 203; A = "vocab:1182"; B = "SARS-CoV-2"; D = "vocab:1178"; F = "http://purl.org/zonmw/covid19/10269" },
    @{ Row = 204; A = "vocab:1183"; B = "SARS-CoV-1"; D = "vocab:1178" },
    @{ Row = 205; A = "vocab:1184"; B = "MERS-CoV"; D = "vocab:1178" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    if ($r.ContainsKey("D")) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    if ($r.ContainsKey("F")) {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
}
